$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update code coverage values in column G
$ws.Range("G3").Value = 0.075
$ws.Range("G4").Value = 0.896
$ws.Range("G6").Value = 0.298
$ws.Range("G9").Value = 0.726

# Update the active cell selection from G11 to G10
$ws.Range("G10").Select()
